$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.2050900623679297
$ws.Range("D2").Value = 0.0374661289536391
$ws.Range("E2").Value = 0.1056517465228595
$ws.Range("F2").Value = 3.506972899598253
$ws.Range("G2").Value = 3.009965487123139
$ws.Range("H2").Value = 2.350412002670367
$ws.Range("J2").Value = 0.2083103321579927
$ws.Range("L2").Value = 0.1919518520350607
$ws.Range("M2").Value = 11.74440127962214

$ws.Range("C3").Value = 0.2126857326095148
$ws.Range("D3").Value = 0.03347376848962824
$ws.Range("E3").Value = 0.1016751828309026
$ws.Range("F3").Value = 3.612214181360372
$ws.Range("G3").Value = 3.094436563046983
$ws.Range("H3").Value = 2.411908727009376
$ws.Range("J3").Value = 0.2055319741289452
$ws.Range("L3").Value = 0.1724249042778325
$ws.Range("M3").Value = 10.48162408359542

$ws.Range("C4").Value = 0.217589545739024
$ws.Range("D4").Value = 0.03101279777876442
$ws.Range("E4").Value = 0.09930655409419131
$ws.Range("F4").Value = 3.681911354656293
$ws.Range("G4").Value = 3.15124241857319
$ws.Range("H4").Value = 2.452604518884513
$ws.Range("J4").Value = 0.2041113387331848
$ws.Range("L4").Value = 0.1604613179831773
$ws.Range("M4").Value = 9.703913012778685

$ws.Range("C5").Value = 0.2196477257750722
$ws.Range("D5").Value = 0.03000734246437986
$ws.Range("E5").Value = 0.09835927004155565
$ws.Range("F5").Value = 3.711573909939645
$ws.Range("G5").Value = 3.17561465223018
$ws.Range("H5").Value = 2.469918330776011
$ws.Range("J5").Value = 0.2036024896885351
$ws.Range("L5").Value = 0.1555918902084557
$ws.Range("M5").Value = 9.386351817242883

$ws.Range("C6").Value = 0.2199930838801158
$ws.Range("D6").Value = 0.02984022629672722
$ws.Range("E6").Value = 0.09820304726708429
$ws.Range("F6").Value = 3.716574997291119
$ws.Range("G6").Value = 3.179734986700581
$ws.Range("H6").Value = 2.472837117911226
$ws.Range("J6").Value = 0.2035221808420005
$ws.Range("L6").Value = 0.1547836581691087
$ws.Range("M6").Value = 9.333581314464823

$ws.Range("C7").Value = 0.2176170614802793
$ws.Range("D7").Value = 0.03099924855482072
$ws.Range("E7").Value = 0.09929370651608593
$ws.Range("F7").Value = 3.682306312749915
$ws.Range("G7").Value = 3.151566182273001
$ws.Range("H7").Value = 2.452835074553604
$ws.Range("J7").Value = 0.2041041946349509
$ws.Range("L7").Value = 0.1603956244594968
$ws.Range("M7").Value = 9.699632908016667

$ws.Range("C8").Value = 0.207658911217024
$ws.Range("D8").Value = 0.03609146266961716
$ws.Range("E8").Value = 0.1042652218304525
$ws.Range("F8").Value = 3.542195926231472
$ws.Range("G8").Value = 3.038054300025465
$ws.Range("H8").Value = 2.371001454607352
$ws.Range("J8").Value = 0.2072921100795782
$ws.Range("L8").Value = 0.1852130944790531
$ws.Range("M8").Value = 11.30945743798952

$ws.Range("C9").Value = 0.1900593136582742
$ws.Range("D9").Value = 0.04600895169735963
$ws.Range("E9").Value = 0.1146140577567394
$ws.Range("F9").Value = 3.308494683596734
$ws.Range("G9").Value = 2.855519085681635
$ws.Range("H9").Value = 2.234211633806439
$ws.Range("J9").Value = 0.2158876256291791
$ws.Range("L9").Value = 0.2341254932093477
$ws.Range("M9").Value = 14.44995130438139

$ws.Range("C10").Value = 0.1783374786089951
$ws.Range("D10").Value = 0.05326620583808506
$ws.Range("E10").Value = 0.1226147641347168
$ws.Range("F10").Value = 3.162886337413639
$ws.Range("G10").Value = 2.747046645138596
$ws.Range("H10").Value = 2.148691919974738
$ws.Range("J10").Value = 0.2237511272690256
$ws.Range("L10").Value = 0.2702744150277852
$ws.Range("M10").Value = 16.75113314888813

$ws.Range("C11").Value = 0.1732748314417627
$ws.Range("D11").Value = 0.056564321139291
$ws.Range("E11").Value = 0.1263485440614005
$ws.Range("F11").Value = 3.10254962392915
$ws.Range("G11").Value = 2.703538948689555
$ws.Range("H11").Value = 2.113160708408998
$ws.Range("J11").Value = 0.2276929129800749
$ws.Range("L11").Value = 0.2867815857658513
$ws.Range("M11").Value = 17.79762128552102

$ws.Range("C12").Value = 0.171397185341867
$ws.Range("D12").Value = 0.05781300445003978
$ws.Range("E12").Value = 0.1277766589060434
$ws.Range("F12").Value = 3.080570926390806
$ws.Range("G12").Value = 2.68792655599492
$ws.Range("H12").Value = 2.10020137011989
$ws.Range("J12").Value = 0.2292404888941775
$ws.Range("L12").Value = 0.2930428251684702
$ws.Range("M12").Value = 18.19393332172814

$ws.Range("C13").Value = 0.1717998033736166
$ws.Range("D13").Value = 0.05754408518109244
$ws.Range("E13").Value = 0.1274684455602895
$ws.Range("F13").Value = 3.085265420067699
$ws.Range("G13").Value = 2.691250195042102
$ws.Range("H13").Value = 2.102970181872394
$ws.Range("J13").Value = 0.2289047077454995
$ws.Range("L13").Value = 0.2916938738903525
$ws.Range("M13").Value = 18.10857781323119

$ws.Range("C14").Value = 0.17311956131541
$ws.Range("D14").Value = 0.05666705457126398
$ws.Range("E14").Value = 0.1264657470875576
$ws.Range("F14").Value = 3.100723874814449
$ws.Range("G14").Value = 2.702237072934196
$ws.Range("H14").Value = 2.112084541835742
$ws.Range("J14").Value = 0.2278191178680089
$ws.Range("L14").Value = 0.2872964868078896
$ws.Range("M14").Value = 17.83022505416807

$ws.Range("C15").Value = 0.1739331120608512
$ws.Range("D15").Value = 0.05612982345265038
$ws.Range("E15").Value = 0.1258534369026876
$ws.Range("F15").Value = 3.110306492698925
$ws.Range("G15").Value = 2.709079960983473
$ws.Range("H15").Value = 2.117732209237914
$ws.Range("J15").Value = 0.2271613891904707
$ws.Range("L15").Value = 0.2846043461485408
$ws.Range("M15").Value = 17.65973222341489

$ws.Range("C16").Value = 0.178673819375188
$ws.Range("D16").Value = 0.05305061993514926
$ws.Range("E16").Value = 0.1223727044737473
$ws.Range("F16").Value = 3.166950003738862
$ws.Range("G16").Value = 2.750009436149185
$ws.Range("H16").Value = 2.151082726672087
$ws.Range("J16").Value = 0.2235010677541567
$ws.Range("L16").Value = 0.2691970053630257
$ws.Range("M16").Value = 16.68274248725231

$ws.Range("C17").Value = 0.1816516419625991
$ws.Range("D17").Value = 0.051160953710081
$ws.Range("E17").Value = 0.1202619771590463
$ws.Range("F17").Value = 3.203225078274414
$ws.Range("G17").Value = 2.776629929482681
$ws.Range("H17").Value = 2.172413341881935
$ws.Range("J17").Value = 0.2213506846093765
$ws.Range("L17").Value = 0.2597620832743814
$ws.Range("M17").Value = 16.083359111444

$ws.Range("C18").Value = 0.1833897469340826
$ws.Range("D18").Value = 0.05007374797528996
$ws.Range("E18").Value = 0.1190567681842012
$ws.Range("F18").Value = 3.22464370429239
$ws.Range("G18").Value = 2.792489272936422
$ws.Range("H18").Value = 2.184998891306691
$ws.Range("J18").Value = 0.220148005881569
$ws.Range("L18").Value = 0.2543411966325948
$ws.Range("M18").Value = 15.73857347315288

$ws.Range("C19").Value = 0.1839825710417848
$ws.Range("D19").Value = 0.0497055772033832
$ws.Range("E19").Value = 0.1186502014678439
$ws.Range("F19").Value = 3.231990235015104
$ws.Range("G19").Value = 2.797952374625424
$ws.Range("H19").Value = 2.189314224288495
$ws.Range("J19").Value = 0.2197465944380355
$ws.Range("L19").Value = 0.2525067413941144
$ws.Range("M19").Value = 15.62182651725129

$ws.Range("C20").Value = 0.1813320201455859
$ws.Range("D20").Value = 0.05136214339782441
$ws.Range("E20").Value = 0.1204857496135929
$ws.Range("F20").Value = 3.199306031679356
$ws.Range("G20").Value = 2.77373925384677
$ws.Range("H20").Value = 2.170109803681925
$ws.Range("J20").Value = 0.2215760432034699
$ws.Range("L20").Value = 0.2607658343249284
$ws.Range("M20").Value = 16.14716769943408

$ws.Range("C21").Value = 0.1727308396150633
$ws.Range("D21").Value = 0.05692466441909971
$ws.Range("E21").Value = 0.1267598726412018
$ws.Range("F21").Value = 3.096159583559853
$ws.Range("G21").Value = 2.698986345493466
$ws.Range("H21").Value = 2.109393897313737
$ws.Range("J21").Value = 0.228136471311899
$ws.Range("L21").Value = 0.2885878147012022
$ws.Range("M21").Value = 17.91198255305932

$ws.Range("C22").Value = 0.1673398639472481
$ws.Range("D22").Value = 0.06055879945405707
$ws.Range("E22").Value = 0.1309435656546043
$ws.Range("F22").Value = 3.033825786424018
$ws.Range("G22").Value = 2.6551729611366
$ws.Range("H22").Value = 2.072606252795282
$ws.Range("J22").Value = 0.2327453873880501
$ws.Range("L22").Value = 0.3068319924022092
$ws.Range("M22").Value = 19.06560467370832

$ws.Range("C23").Value = 0.1701958129248951
$ws.Range("D23").Value = 0.05861923410014924
$ws.Range("E23").Value = 0.1287028053311445
$ws.Range("F23").Value = 3.066622739519445
$ws.Range("G23").Value = 2.678087612445012
$ws.Range("H23").Value = 2.091972122301001
$ws.Range("J23").Value = 0.2302552671590661
$ws.Range("L23").Value = 0.2970887012894536
$ws.Range("M23").Value = 18.4498480588046

$ws.Range("C24").Value = 0.1814764398780646
$ws.Range("D24").Value = 0.05127118808395892
$ws.Range("E24").Value = 0.1203845563578341
$ws.Range("F24").Value = 3.201076078723531
$ws.Range("G24").Value = 2.775044401796691
$ws.Range("H24").Value = 2.17115023022879
$ws.Range("J24").Value = 0.2214740540958502
$ws.Range("L24").Value = 0.2603120280566316
$ws.Range("M24").Value = 16.11832044107683

$ws.Range("C25").Value = 0.1946106935695582
$ws.Range("D25").Value = 0.04333206439156356
$ws.Range("E25").Value = 0.1117470877665028
$ws.Range("F25").Value = 3.367213810894
$ws.Range("G25").Value = 2.90048760628855
$ws.Range("H25").Value = 2.268626304455552
$ws.Range("J25").Value = 0.2132994851053667
$ws.Range("L25").Value = 0.2208609342761889
$ws.Range("M25").Value = 13.60174103047189

